$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-23 Thursday" "2025-10-24 Friday"

Replace-Text "703÷5=" "672÷8="
Replace-Text "219÷3=" "885÷3="
Replace-Text "885÷5=" "759÷6="
Replace-Text "109÷8=" "952÷7="
Replace-Text "686÷3=" "961÷2="
Replace-Text "941÷2=" "449÷4="
Replace-Text "637÷7=" "561÷5="
Replace-Text "122÷6=" "893÷9="
Replace-Text "808÷7=" "937÷5="
Replace-Text "900÷6=" "661÷3="
Replace-Text "762÷9=" "676÷3="
Replace-Text "499÷4=" "660÷5="
Replace-Text "501÷7=" "436÷8="
Replace-Text "309÷4=" "931÷4="
Replace-Text "361÷8=" "984÷7="
Replace-Text "847÷3=" "288÷6="
Replace-Text "860÷9=" "873÷8="
Replace-Text "398÷7=" "965÷2="
Replace-Text "742÷4=" "101÷6="
Replace-Text "869÷2=" "334÷6="
Replace-Text "190÷6=" "457÷3="
Replace-Text "228÷8=" "245÷5="
Replace-Text "734÷2=" "353÷4="
Replace-Text "591÷9=" "793÷3="
Replace-Text "239÷8=" "452÷7="
